$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Scene" data table lost two rows ("CloneScene"/Scene2 and
# "RebellerNoob"/SelectScene) and two of the remaining rows got their
# data values updated for the server config.

# Remove the CloneScene/Scene2 row (original row 2).
$ws.Rows.Item(2).Delete()

# Remove the RebellerNoob/SelectScene row (now row 3, after the row-2
# delete shifted everything up by one).
$ws.Rows.Item(3).Delete()

# Row 2 is now the PioneerNoob/villageScene entry: update its RelivePos
# value and drop the Text-format styling that F:G used to carry.
$ws.Range("E2").Value = "20,0,60"
$ws.Range("F2:G2").ClearFormats()

# Row 3 is now the old Demo1 entry: its ID changes from 4 to 2.
$ws.Range("B3").Value = "2"
$ws.Range("F3").Value = "Demo1"
$ws.Range("G3").Value = "Demo1"

# Match the saved selection state.
[void]$ws.Range("F5").Select()
